# Apply cryptos list update (Thu Oct 26 03:30:44 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "34.725.52"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +2.11%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.801.00"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("E4").Value = "  -0.23%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "224.80"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -1.39%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.556"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("E7").Value = "  -0.30%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "32.52"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +4.66%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.288"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +2.85%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0722"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +9.43%  "

$ws.Range("E11").Value = "  +0.68%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "2.055.27"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +0.29%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "11.18"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "1.788.27"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -0.10%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.641"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +1.33%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "34.740.52"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +2.03%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "4.33"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +2.99%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "69.52"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +0.26%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "254.98"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +0.94%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0814"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +9.42%  "

$ws.Range("E21").Value = "  -0.10%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.85"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +4.26%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "4.25"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("E24").Value = "  +0.31%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "160.84"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +2.16%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "16.51"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -0.40%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "7.16"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +2.24%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.115"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("E29").Value = "  -0.32%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.0537"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +3.96%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "3.82"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -1.47%  "

$ws.Range("E32").Value = "  -0.59%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.63"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +0.32%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.89"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +2.66%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.445.95"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -1.90%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.06"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("E37").Value = "  +3.34%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.641"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +1.46%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "85.47"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +2.07%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.80"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -1.08%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.943"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +4.79%  "

$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "2.33"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.13"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +3.92%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "6.01"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +5.56%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.06"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0495"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -3.97%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.953.49"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "106.30"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +9.09%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "12.06"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +2.76%  "

$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0₆0127"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +9.76%  "
